$wb = $excel.ActiveWorkbook

# xlPasteFormats constant used below for style-only paste.
$xlPasteFormats = -4122

# --- 1. Insert a new "2022-Q1" sheet right after "2021-Q4", before "总计" ---
$ws2021Q4 = $wb.Worksheets.Item("2021-Q4")
$wsTotal = $wb.Worksheets.Item("总计")

$wsNew = $wb.Worksheets.Add($wsTotal)
$wsNew.Name = "2022-Q1"

# --- 2. Populate the new "2022-Q1" sheet (same column layout as "2021-Q4") ---
$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

$wsNew.Range("A2").Value = 0
$wsNew.Range("B2").Value = "'011444"
$wsNew.Range("C2").Value = "创金合信瑞裕混合A"
$wsNew.Range("D2").Value = "'0.03"
$wsNew.Range("E2").Value = "'68.73"
$wsNew.Range("F2").Value = "'2.27"
$wsNew.Range("G2").Value = "'0.0007"
$wsNew.Range("H2").Value = 9

$wsNew.Range("A3").Value = 1
$wsNew.Range("B3").Value = "'011445"
$wsNew.Range("C3").Value = "创金合信瑞裕混合C"
$wsNew.Range("D3").Value = "'0.00"
$wsNew.Range("E3").Value = "'68.73"
$wsNew.Range("F3").Value = "'2.27"
$wsNew.Range("G3").Value = 0
$wsNew.Range("H3").Value = 9

# Drop the auto "quote-prefix" number format that typing a leading apostrophe
# triggers, then re-apply the real header / index-column look from the
# "2021-Q4" sheet so the new tab matches the rest of the workbook.
# (Deliberately skip A1: it is never used/populated, same as on "2021-Q4".)
$wsNew.Range("B1:H3").ClearFormats()
$wsNew.Range("A2:A3").ClearFormats()

$ws2021Q4.Range("B1:H1").Copy()
$wsNew.Range("B1:H1").PasteSpecial($xlPasteFormats)

$ws2021Q4.Range("A2").Copy()
$wsNew.Range("A2:A3").PasteSpecial($xlPasteFormats)

# --- 3. Update "总计" sheet: insert a new row for 2022-Q1 above the 2021-Q4 row ---
# Re-resolve by name: worksheet handles track by index, and the sheet insert
# above shifted "总计" from index 2 to index 3.
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0
$wsTotal.Range("A2:D2").ClearFormats()

$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial($xlPasteFormats)

$wsTotal.Range("A3").Value = 1

$excel.CutCopyMode = $false
